$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.765.52'
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = '2.286.69'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.619'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.02'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0944'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.56%  '
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").Value = '2.630.71'
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.872'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.72%  '
$ws.Range("D17").Value = '2.288.55'
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '43.618.12'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("E19").Value = '  -1.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +11.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.53'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +12.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.97'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.56%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.02%  '
$ws.Range("E28").Value = '  +18.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '41.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +10.28%  '
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("E31").Value = '  +1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0935'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.57'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E36").Value = '  +0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  +1.94%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.107'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("E40").Value = '  +8.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.50'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +18.74%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.96%  '
$ws.Range("E43").Value = '  +3.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.242'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.10%  '
$ws.Range("E45").Value = '  +21.70%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.26'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.95%  '
